$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns B:G, rows 2-11 (regenerated s_val data filtering save games)
$data = @(
    @(1.445647641019636, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 0, 3.755628166162433),
    @(0.04172184405617529, 0.04103571897497393, 3.223369029078222, 0.5333859586016987, 1, 3.83951255071107),
    @(3.272327238179451, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 1, 5.582307763322248),
    @(0.2881169905109251, 0.3048912486333797, 0.1496068669990043, 0.5333859586016987, 0, 1.276001064745008),
    @(3.272327238179451, 1.626987699542094, 3.223369029078222, 0.5333859586016987, 1, 8.656069925401464),
    @(3.272327238179451, 1.626987699542094, 3.223369029078222, 0.5333859586016987, 1, 8.656069925401464),
    @(3.272327238179451, 1.626987699542094, 3.223369029078222, 13.86384647080068, 1, 21.98653043760045),
    @(1.445647641019636, 1.626987699542094, 3993.344853322108, 13.86384647080068, 1, 4010.28133513347),
    @(1.445647641019636, 1.626987699542094, 3.223369029078222, 13.86384647080068, 1, 20.15985084044064),
    @(0.6545652718822623, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 1, 2.964545797025059)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
}
